$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 26, pushing existing rows 26-34 down to 27-35.
$ws.Rows.Item(26).Insert()

# Populate the newly inserted row 26 with its data.
$ws.Cells.Item(26, 1).Value = 6
$ws.Cells.Item(26, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(26, 3).Value = "Metropolitana"
$ws.Cells.Item(26, 4).Value = 44784
$ws.Cells.Item(26, 5).Value = 13
$ws.Cells.Item(26, 6).Value = 100112035
$ws.Cells.Item(26, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(26, 8).Value = "Sin especificar"
$ws.Cells.Item(26, 9).Value = "Primera"
$ws.Cells.Item(26, 10).Value = 220
$ws.Cells.Item(26, 11).Value = 17000
$ws.Cells.Item(26, 12).Value = 18000
$ws.Cells.Item(26, 13).Value = 17455
$ws.Cells.Item(26, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(26, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(26, 16).Value = 1164
$ws.Cells.Item(26, 17).Value = 15
$ws.Cells.Item(26, 18).Value = "Hortaliza"
